# Update "想去人数" (interest count) values in column F across the
# "展览", "演出" and "全部类型" sheets to match the refreshed data pull.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws2 = $wb.Worksheets.Item("演出")
$ws4 = $wb.Worksheets.Item("全部类型")

# 展览 sheet (sheet1)
$sheet1Updates = @{
    3  = 1424
    4  = 20037
    5  = 798
    6  = 311
    7  = 1097
    9  = 7564
    10 = 512
    11 = 736
    12 = 264
    14 = 159
    16 = 11
    18 = 192
    19 = 1335
    20 = 419
    22 = 679
    24 = 67
    27 = 1100
    28 = 31
    30 = 182
    32 = 559
    33 = 64
    34 = 2841
    38 = 12628
    39 = 1332
    40 = 80
    41 = 25
    42 = 54
    43 = 262
    45 = 3998
    46 = 320
}

foreach ($row in $sheet1Updates.Keys) {
    $ws1.Range("F$row").Value = $sheet1Updates[$row]
}

# 演出 sheet (sheet2)
$ws2.Range("F2").Value = 166

# 全部类型 sheet (sheet4)
$sheet4Updates = @{
    3  = 1424
    4  = 20037
    5  = 798
    6  = 311
    7  = 1097
    9  = 7564
    10 = 512
    11 = 736
    12 = 264
    14 = 159
    16 = 11
    18 = 192
    19 = 1335
    20 = 419
    22 = 679
    24 = 67
    27 = 1100
    28 = 31
    30 = 182
    31 = 166
    32 = 559
    34 = 64
    36 = 2841
    37 = 25
    40 = 12628
    41 = 1332
    42 = 80
    43 = 25
    44 = 54
    45 = 262
    47 = 3998
    48 = 320
}

foreach ($row in $sheet4Updates.Keys) {
    $ws4.Range("F$row").Value = $sheet4Updates[$row]
}
